$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at positions 12 and 13 (pushing the old "Programa resumido:" block,
# and everything below it, down by two rows).
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(13).Insert()

# Row 12 only carries the new label in column A ("Docentes responsáveis:").
$ws.Range("A12").Value2 = "Docentes responsáveis:"
$ws.Range("B12").Clear()
$ws.Range("C12").Clear()

# Row 13 only carries the value (same text repeated in the PT and EN columns).
$ws.Range("A13").Clear()
$ws.Range("B13").Value2 = "210064 - Eduardo Rezende Triboni"
$ws.Range("C13").Value2 = "210064 - Eduardo Rezende Triboni"
